# Add a new "2020" column (Q) to the summary table, mirroring the
# formatting already used by the adjacent "2019" column (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column P (the existing last data column) into
# the new column Q for every row of the table (header row 4 through the
# last data row 14).
$ws.Range("P4:P14").Copy()
$ws.Range("Q4:Q14").PasteSpecial(-4122)

# Year header
$ws.Range("Q4").Value = 2020

# Data values for the new 2020 column
$ws.Range("Q5").Value = 99.3
$ws.Range("Q6").Value = 99.371420589467803
$ws.Range("Q7").Value = 99.319469393395053
$ws.Range("Q8").Value = 99.442213297634979
$ws.Range("Q9").Value = 98.766881972988841
$ws.Range("Q10").Value = 99.212798374809537
$ws.Range("Q11").Value = 99.799160124155549
$ws.Range("Q12").Value = 99.146991622239156
$ws.Range("Q13").Value = 99.538370126605429
$ws.Range("Q14").Value = 99.765563948945029

# Match the selection state recorded in the saved file
$ws.Range("P7").Select()
